$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 8
$ws.Range("G2").Value = 0.444
$ws.Range("H2").Value = 48.6
$ws.Range("J2").Value = 77.40000000000001
$ws.Range("K2").Value = 0.441
$ws.Range("M2").Value = 13.7
$ws.Range("O2").Value = 20.9
$ws.Range("P2").Value = 27.1
$ws.Range("Q2").Value = 0.77
$ws.Range("R2").Value = 11.8
$ws.Range("S2").Value = 28.6
$ws.Range("T2").Value = 40.4
$ws.Range("U2").Value = 19.4
$ws.Range("V2").Value = 15.7
$ws.Range("W2").Value = 7.7
$ws.Range("X2").Value = 6.1
$ws.Range("Y2").Value = 5.9
$ws.Range("Z2").Value = 19.8
$ws.Range("AA2").Value = 21.7
$ws.Range("AB2").Value = 93.3
$ws.Range("AC2").Value = -1.8
$ws.Range("AD2").Value = 22
$ws.Range("AE2").Value = 18
$ws.Range("AG2").Value = 18
$ws.Range("AJ2").Value = 27
$ws.Range("AK2").Value = 21
$ws.Range("AP2").Value = 10
$ws.Range("AQ2").Value = 9
$ws.Range("AR2").Value = 13
$ws.Range("AT2").Value = 26
$ws.Range("AV2").Value = 17
$ws.Range("AW2").Value = 13
$ws.Range("AX2").Value = 1
$ws.Range("AY2").Value = 28
$ws.Range("AZ2").Value = 4
$ws.Range("BA2").Value = 16
$ws.Range("BB2").Value = 23
$ws.Range("BF2").Value = "2007-12-08"
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = 16
$ws.Range("G3").Value = 0.889
$ws.Range("H3").Value = 48.6
$ws.Range("I3").Value = 36.6
$ws.Range("J3").Value = 75.8
$ws.Range("K3").Value = 0.482
$ws.Range("M3").Value = 18.9
$ws.Range("N3").Value = 0.394
$ws.Range("O3").Value = 21.7
$ws.Range("P3").Value = 28.3
$ws.Range("Q3").Value = 0.766
$ws.Range("R3").Value = 9.1
$ws.Range("S3").Value = 33.4
$ws.Range("T3").Value = 42.5
$ws.Range("U3").Value = 23.8
$ws.Range("V3").Value = 15.9
$ws.Range("W3").Value = 9.6
$ws.Range("X3").Value = 4.3
$ws.Range("Y3").Value = 5.1
$ws.Range("Z3").Value = 21.9
$ws.Range("AA3").Value = 23.3
$ws.Range("AB3").Value = 102.2
$ws.Range("AC3").Value = 14.1
$ws.Range("AD3").Value = 22
$ws.Range("AE3").Value = 2
$ws.Range("AI3").Value = 16
$ws.Range("AM3").Value = 13
$ws.Range("AP3").Value = 8
$ws.Range("AQ3").Value = 12
$ws.Range("AS3").Value = 6
$ws.Range("AT3").Value = 14
$ws.Range("AU3").Value = 4
$ws.Range("AV3").Value = 21
$ws.Range("AX3").Value = 24
$ws.Range("AY3").Value = 16
$ws.Range("BB3").Value = 10
$ws.Range("BF3").Value = "2007-12-08"
$ws.Range("D4").Value = 17
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 0.353
$ws.Range("I4").Value = 34.1
$ws.Range("J4").Value = 79.09999999999999
$ws.Range("K4").Value = 0.431
$ws.Range("L4").Value = 6.4
$ws.Range("M4").Value = 17.8
$ws.Range("N4").Value = 0.361
$ws.Range("O4").Value = 17.2
$ws.Range("P4").Value = 25.4
$ws.Range("Q4").Value = 0.677
$ws.Range("R4").Value = 12.6
$ws.Range("S4").Value = 28.4
$ws.Range("T4").Value = 41
$ws.Range("U4").Value = 19.1
$ws.Range("V4").Value = 16.6
$ws.Range("W4").Value = 7.8
$ws.Range("X4").Value = 4.7
$ws.Range("AA4").Value = 21.1
$ws.Range("AB4").Value = 91.8
$ws.Range("AC4").Value = -5.9
$ws.Range("AD4").Value = 28
$ws.Range("AE4").Value = 23
$ws.Range("AG4").Value = 23
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 27
$ws.Range("AM4").Value = 14
$ws.Range("AO4").Value = 21
$ws.Range("AP4").Value = 18
$ws.Range("AR4").Value = 4
$ws.Range("AT4").Value = 22
$ws.Range("AX4").Value = 20
$ws.Range("BA4").Value = 20
$ws.Range("BC4").Value = 27
$ws.Range("BF4").Value = "2007-12-08"
$ws.Range("D5").Value = 17
$ws.Range("F5").Value = 11
$ws.Range("G5").Value = 0.353
$ws.Range("I5").Value = 33.8
$ws.Range("J5").Value = 84.59999999999999
$ws.Range("K5").Value = 0.4
$ws.Range("L5").Value = 5.2
$ws.Range("N5").Value = 0.319
$ws.Range("Q5").Value = 0.746
$ws.Range("R5").Value = 14.4
$ws.Range("S5").Value = 30.4
$ws.Range("T5").Value = 44.8
$ws.Range("V5").Value = 15.8
$ws.Range("X5").Value = 4.5
$ws.Range("Y5").Value = 5.8
$ws.Range("Z5").Value = 22.1
$ws.Range("AA5").Value = 21.8
$ws.Range("AB5").Value = 89.8
$ws.Range("AC5").Value = -4.8
$ws.Range("AD5").Value = 28
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 18
$ws.Range("AG5").Value = 23
$ws.Range("AL5").Value = 23
$ws.Range("AN5").Value = 27
$ws.Range("AO5").Value = 25
$ws.Range("AP5").Value = 23
$ws.Range("AS5").Value = 17
$ws.Range("AT5").Value = 3
$ws.Range("AV5").Value = 19
$ws.Range("AW5").Value = 10
$ws.Range("AX5").Value = 22
$ws.Range("AY5").Value = 27
$ws.Range("AZ5").Value = 18
$ws.Range("BC5").Value = 20
$ws.Range("BF5").Value = "2007-12-08"
$ws.Range("D6").Value = 20
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 0.45
$ws.Range("I6").Value = 35.7
$ws.Range("J6").Value = 82.2
$ws.Range("K6").Value = 0.434
$ws.Range("L6").Value = 7.4
$ws.Range("M6").Value = 19.8
$ws.Range("N6").Value = 0.371
$ws.Range("O6").Value = 17.5
$ws.Range("P6").Value = 24.7
$ws.Range("Q6").Value = 0.706
$ws.Range("S6").Value = 31.2
$ws.Range("T6").Value = 43.4
$ws.Range("U6").Value = 18.9
$ws.Range("V6").Value = 15.2
$ws.Range("Y6").Value = 5.4
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 19.3
$ws.Range("AB6").Value = 96.09999999999999
$ws.Range("AC6").Value = -5.3
$ws.Range("AD6").Value = 2
$ws.Range("AF6").Value = 18
$ws.Range("AG6").Value = 16
$ws.Range("AH6").Value = 10
$ws.Range("AI6").Value = 20
$ws.Range("AK6").Value = 25
$ws.Range("AM6").Value = 8
$ws.Range("AN6").Value = 8
$ws.Range("AO6").Value = 18
$ws.Range("AP6").Value = 20
$ws.Range("AR6").Value = 10
$ws.Range("AS6").Value = 13
$ws.Range("AV6").Value = 12
$ws.Range("AX6").Value = 16
$ws.Range("AZ6").Value = 17
$ws.Range("BB6").Value = 18
$ws.Range("BC6").Value = 24
$ws.Range("BF6").Value = "2007-12-08"
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 12
$ws.Range("G7").Value = 0.6
$ws.Range("H7").Value = 48.3
$ws.Range("I7").Value = 36.6
$ws.Range("J7").Value = 79.59999999999999
$ws.Range("K7").Value = 0.46
$ws.Range("L7").Value = 5.7
$ws.Range("M7").Value = 17.1
$ws.Range("N7").Value = 0.33
$ws.Range("O7").Value = 23.7
$ws.Range("P7").Value = 28
$ws.Range("Q7").Value = 0.846
$ws.Range("R7").Value = 10.3
$ws.Range("S7").Value = 32.5
$ws.Range("T7").Value = 42.8
$ws.Range("U7").Value = 20.4
$ws.Range("W7").Value = 5.6
$ws.Range("Y7").Value = 4.3
$ws.Range("Z7").Value = 23.3
$ws.Range("AA7").Value = 22.7
$ws.Range("AB7").Value = 102.5
$ws.Range("AC7").Value = 3.3
$ws.Range("AD7").Value = 2
$ws.Range("AE7").Value = 8
$ws.Range("AF7").Value = 8
$ws.Range("AG7").Value = 8
$ws.Range("AH7").Value = 16
$ws.Range("AI7").Value = 15
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 9
$ws.Range("AL7").Value = 20
$ws.Range("AO7").Value = 2
$ws.Range("AP7").Value = 9
$ws.Range("AR7").Value = 20
$ws.Range("AS7").Value = 9
$ws.Range("AT7").Value = 12
$ws.Range("AU7").Value = 17
$ws.Range("AX7").Value = 12
$ws.Range("AY7").Value = 8
$ws.Range("AZ7").Value = 24
$ws.Range("BA7").Value = 11
$ws.Range("BB7").Value = 9
$ws.Range("BF7").Value = "2007-12-08"
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 12
$ws.Range("G8").Value = 0.6
$ws.Range("I8").Value = 38.8
$ws.Range("J8").Value = 84.5
$ws.Range("L8").Value = 6.7
$ws.Range("M8").Value = 19.2
$ws.Range("N8").Value = 0.349
$ws.Range("O8").Value = 23.2
$ws.Range("P8").Value = 31.5
$ws.Range("Q8").Value = 0.737
$ws.Range("R8").Value = 10.9
$ws.Range("U8").Value = 23.6
$ws.Range("V8").Value = 16.4
$ws.Range("W8").Value = 10.2
$ws.Range("X8").Value = 6
$ws.Range("Y8").Value = 4.6
$ws.Range("Z8").Value = 22.8
$ws.Range("AA8").Value = 25
$ws.Range("AB8").Value = 107.4
$ws.Range("AD8").Value = 2
$ws.Range("AE8").Value = 8
$ws.Range("AF8").Value = 8
$ws.Range("AG8").Value = 8
$ws.Range("AM8").Value = 11
$ws.Range("AN8").Value = 18
$ws.Range("AR8").Value = 18
$ws.Range("AS8").Value = 5
$ws.Range("AT8").Value = 6
$ws.Range("AU8").Value = 5
$ws.Range("AV8").Value = 24
$ws.Range("AX8").Value = 2
$ws.Range("AY8").Value = 13
$ws.Range("AZ8").Value = 22
$ws.Range("BF8").Value = "2007-12-08"
$ws.Range("AD9").Value = 14
$ws.Range("AJ9").Value = 15
$ws.Range("AL9").Value = 18
$ws.Range("AM9").Value = 21
$ws.Range("AN9").Value = 11
$ws.Range("AP9").Value = 19
$ws.Range("AS9").Value = 27
$ws.Range("AX9").Value = 3
$ws.Range("BF9").Value = "2007-12-08"
$ws.Range("AD10").Value = 14
$ws.Range("AF10").Value = 8
$ws.Range("AH10").Value = 6
$ws.Range("AN10").Value = 12
$ws.Range("AP10").Value = 13
$ws.Range("AS10").Value = 18
$ws.Range("AT10").Value = 10
$ws.Range("AX10").Value = 19
$ws.Range("AY10").Value = 23
$ws.Range("AZ10").Value = 23
$ws.Range("BA10").Value = 10
$ws.Range("BB10").Value = 2
$ws.Range("BF10").Value = "2007-12-08"
$ws.Range("AD11").Value = 2
$ws.Range("AI11").Value = 14
$ws.Range("AK11").Value = 18
$ws.Range("AN11").Value = 28
$ws.Range("AP11").Value = 24
$ws.Range("AS11").Value = 11
$ws.Range("AT11").Value = 4
$ws.Range("AW11").Value = 11
$ws.Range("BB11").Value = 19
$ws.Range("BF11").Value = "2007-12-08"
$ws.Range("AD12").Value = 2
$ws.Range("AN12").Value = 22
$ws.Range("AO12").Value = 13
$ws.Range("AP12").Value = 14
$ws.Range("AR12").Value = 3
$ws.Range("AS12").Value = 4
$ws.Range("AW12").Value = 9
$ws.Range("AY12").Value = 24
$ws.Range("BA12").Value = 9
$ws.Range("BB12").Value = 8
$ws.Range("BF12").Value = "2007-12-08"
$ws.Range("AD13").Value = 22
$ws.Range("AJ13").Value = 25
$ws.Range("AL13").Value = 19
$ws.Range("AM13").Value = 20
$ws.Range("AU13").Value = 22
$ws.Range("AX13").Value = 4
$ws.Range("AY13").Value = 19
$ws.Range("AZ13").Value = 16
$ws.Range("BC13").Value = 21
$ws.Range("BF13").Value = "2007-12-08"
$ws.Range("AD14").Value = 14
$ws.Range("AF14").Value = 8
$ws.Range("AN14").Value = 10
$ws.Range("AR14").Value = 17
$ws.Range("AT14").Value = 5
$ws.Range("AV14").Value = 29
$ws.Range("AX14").Value = 17
$ws.Range("AY14").Value = 15
$ws.Range("AZ14").Value = 19
$ws.Range("BF14").Value = "2007-12-08"
$ws.Range("D15").Value = 19
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 0.316
$ws.Range("I15").Value = 37.7
$ws.Range("J15").Value = 80.59999999999999
$ws.Range("K15").Value = 0.468
$ws.Range("L15").Value = 8.6
$ws.Range("M15").Value = 21.5
$ws.Range("N15").Value = 0.402
$ws.Range("O15").Value = 20.4
$ws.Range("P15").Value = 26.8
$ws.Range("Q15").Value = 0.76
$ws.Range("T15").Value = 41.8
$ws.Range("U15").Value = 21.6
$ws.Range("V15").Value = 15.5
$ws.Range("X15").Value = 5.3
$ws.Range("Z15").Value = 19.2
$ws.Range("AA15").Value = 23.1
$ws.Range("AB15").Value = 104.5
$ws.Range("AC15").Value = -1.7
$ws.Range("AD15").Value = 14
$ws.Range("AE15").Value = 23
$ws.Range("AF15").Value = 26
$ws.Range("AG15").Value = 26
$ws.Range("AH15").Value = 6
$ws.Range("AJ15").Value = 17
$ws.Range("AK15").Value = 7
$ws.Range("AL15").Value = 4
$ws.Range("AO15").Value = 11
$ws.Range("AP15").Value = 12
$ws.Range("AQ15").Value = 14
$ws.Range("AS15").Value = 10
$ws.Range("AT15").Value = 18
$ws.Range("AV15").Value = 16
$ws.Range("AX15").Value = 9
$ws.Range("AY15").Value = 18
$ws.Range("AZ15").Value = 3
$ws.Range("BA15").Value = 8
$ws.Range("BB15").Value = 7
$ws.Range("BC15").Value = 16
$ws.Range("BF15").Value = "2007-12-08"
$ws.Range("AD16").Value = 14
$ws.Range("AI16").Value = 25
$ws.Range("AN16").Value = 20
$ws.Range("AO16").Value = 23
$ws.Range("AP16").Value = 17
$ws.Range("AS16").Value = 19
$ws.Range("AW16").Value = 14
$ws.Range("AX16").Value = 10
$ws.Range("BA16").Value = 19
$ws.Range("BF16").Value = "2007-12-08"
$ws.Range("AD17").Value = 22
$ws.Range("AE17").Value = 18
$ws.Range("AO17").Value = 19
$ws.Range("AS17").Value = 26
$ws.Range("AT17").Value = 20
$ws.Range("AV17").Value = 13
$ws.Range("BC17").Value = 22
$ws.Range("BF17").Value = "2007-12-08"
$ws.Range("D18").Value = 17
$ws.Range("E18").Value = 2
$ws.Range("G18").Value = 0.118
$ws.Range("I18").Value = 36.4
$ws.Range("J18").Value = 82.5
$ws.Range("L18").Value = 5.8
$ws.Range("M18").Value = 16.5
$ws.Range("N18").Value = 0.354
$ws.Range("O18").Value = 14.2
$ws.Range("P18").Value = 19.9
$ws.Range("Q18").Value = 0.716
$ws.Range("R18").Value = 12.5
$ws.Range("S18").Value = 28.6
$ws.Range("T18").Value = 41.1
$ws.Range("U18").Value = 18.3
$ws.Range("V18").Value = 16.4
$ws.Range("X18").Value = 4.7
$ws.Range("Y18").Value = 5.1
$ws.Range("Z18").Value = 25.1
$ws.Range("AA18").Value = 18.3
$ws.Range("AB18").Value = 92.90000000000001
$ws.Range("AC18").Value = -9.199999999999999
$ws.Range("AD18").Value = 28
$ws.Range("AI18").Value = 17
$ws.Range("AL18").Value = 17
$ws.Range("AR18").Value = 8
$ws.Range("AS18").Value = 28
$ws.Range("AT18").Value = 21
$ws.Range("AV18").Value = 26
$ws.Range("AW18").Value = 7
$ws.Range("AX18").Value = 20
$ws.Range("AY18").Value = 17
$ws.Range("BB18").Value = 24
$ws.Range("BC18").Value = 30
$ws.Range("BF18").Value = "2007-12-08"
$ws.Range("AD19").Value = 2
$ws.Range("AG19").Value = 16
$ws.Range("AH19").Value = 10
$ws.Range("AM19").Value = 15
$ws.Range("AN19").Value = 26
$ws.Range("AR19").Value = 19
$ws.Range("AS19").Value = 20
$ws.Range("AT19").Value = 25
$ws.Range("AV19").Value = 27
$ws.Range("AX19").Value = 18
$ws.Range("AY19").Value = 8
$ws.Range("BC19").Value = 26
$ws.Range("BF19").Value = "2007-12-08"
$ws.Range("AD20").Value = 2
$ws.Range("AK20").Value = 22
$ws.Range("AM20").Value = 9
$ws.Range("AO20").Value = 26
$ws.Range("AT20").Value = 13
$ws.Range("AW20").Value = 11
$ws.Range("AZ20").Value = 5
$ws.Range("BF20").Value = "2007-12-08"
$ws.Range("D21").Value = 18
$ws.Range("F21").Value = 12
$ws.Range("G21").Value = 0.333
$ws.Range("H21").Value = 48.6
$ws.Range("I21").Value = 34.9
$ws.Range("J21").Value = 80.59999999999999
$ws.Range("K21").Value = 0.433
$ws.Range("M21").Value = 15.8
$ws.Range("N21").Value = 0.327
$ws.Range("O21").Value = 18.8
$ws.Range("P21").Value = 27
$ws.Range("R21").Value = 12.4
$ws.Range("S21").Value = 29.7
$ws.Range("T21").Value = 42.1
$ws.Range("U21").Value = 17.6
$ws.Range("V21").Value = 15.9
$ws.Range("X21").Value = 2.9
$ws.Range("Z21").Value = 22.7
$ws.Range("AB21").Value = 93.90000000000001
$ws.Range("AC21").Value = -7.8
$ws.Range("AD21").Value = 22
$ws.Range("AE21").Value = 23
$ws.Range("AF21").Value = 24
$ws.Range("AG21").Value = 25
$ws.Range("AI21").Value = 24
$ws.Range("AJ21").Value = 16
$ws.Range("AK21").Value = 26
$ws.Range("AL21").Value = 24
$ws.Range("AM21").Value = 22
$ws.Range("AN21").Value = 25
$ws.Range("AO21").Value = 16
$ws.Range("AP21").Value = 11
$ws.Range("AR21").Value = 9
$ws.Range("AS21").Value = 23
$ws.Range("AT21").Value = 16
$ws.Range("AV21").Value = 20
$ws.Range("AY21").Value = 29
$ws.Range("AZ21").Value = 21
$ws.Range("BB21").Value = 22
$ws.Range("BC21").Value = 29
$ws.Range("BF21").Value = "2007-12-08"
$ws.Range("AE22").Value = 2
$ws.Range("AF22").Value = 4
$ws.Range("AG22").Value = 4
$ws.Range("AH22").Value = 13
$ws.Range("AK22").Value = 8
$ws.Range("AR22").Value = 29
$ws.Range("AV22").Value = 15
$ws.Range("AY22").Value = 7
$ws.Range("BA22").Value = 4
$ws.Range("BC22").Value = 6
$ws.Range("BF22").Value = "2007-12-08"
$ws.Range("D23").Value = 19
$ws.Range("E23").Value = 6
$ws.Range("G23").Value = 0.316
$ws.Range("I23").Value = 35.3
$ws.Range("J23").Value = 79.59999999999999
$ws.Range("K23").Value = 0.443
$ws.Range("L23").Value = 4.2
$ws.Range("N23").Value = 0.306
$ws.Range("O23").Value = 17.3
$ws.Range("P23").Value = 24.5
$ws.Range("Q23").Value = 0.706
$ws.Range("R23").Value = 12.8
$ws.Range("S23").Value = 30
$ws.Range("T23").Value = 42.8
$ws.Range("U23").Value = 20.2
$ws.Range("V23").Value = 15.9
$ws.Range("W23").Value = 6.8
$ws.Range("Z23").Value = 21
$ws.Range("AA23").Value = 20.4
$ws.Range("AB23").Value = 92
$ws.Range("AC23").Value = -2.7
$ws.Range("AD23").Value = 14
$ws.Range("AE23").Value = 23
$ws.Range("AF23").Value = 26
$ws.Range("AG23").Value = 26
$ws.Range("AH23").Value = 6
$ws.Range("AI23").Value = 22
$ws.Range("AK23").Value = 19
$ws.Range("AN23").Value = 29
$ws.Range("AO23").Value = 20
$ws.Range("AS23").Value = 21
$ws.Range("AT23").Value = 11
$ws.Range("AV23").Value = 22
$ws.Range("AW23").Value = 20
$ws.Range("AY23").Value = 12
$ws.Range("BB23").Value = 26
$ws.Range("BC23").Value = 18
$ws.Range("BF23").Value = "2007-12-08"
$ws.Range("D24").Value = 20
$ws.Range("F24").Value = 4
$ws.Range("G24").Value = 0.8
$ws.Range("I24").Value = 42.7
$ws.Range("J24").Value = 85.90000000000001
$ws.Range("K24").Value = 0.497
$ws.Range("L24").Value = 8.699999999999999
$ws.Range("M24").Value = 23.3
$ws.Range("N24").Value = 0.373
$ws.Range("O24").Value = 17.1
$ws.Range("P24").Value = 21.8
$ws.Range("Q24").Value = 0.786
$ws.Range("R24").Value = 8.6
$ws.Range("S24").Value = 33.4
$ws.Range("T24").Value = 41.9
$ws.Range("W24").Value = 8.199999999999999
$ws.Range("X24").Value = 5.5
$ws.Range("Y24").Value = 3.8
$ws.Range("Z24").Value = 18.2
$ws.Range("AA24").Value = 19.8
$ws.Range("AB24").Value = 111.1
$ws.Range("AC24").Value = 6.5
$ws.Range("AD24").Value = 2
$ws.Range("AE24").Value = 2
$ws.Range("AK24").Value = 1
$ws.Range("AN24").Value = 6
$ws.Range("AO24").Value = 22
$ws.Range("AP24").Value = 26
$ws.Range("AT24").Value = 17
$ws.Range("AX24").Value = 6
$ws.Range("BB24").Value = 1
$ws.Range("BC24").Value = 5
$ws.Range("BF24").Value = "2007-12-08"
$ws.Range("AD25").Value = 14
$ws.Range("AF25").Value = 24
$ws.Range("AI25").Value = 21
$ws.Range("AJ25").Value = 26
$ws.Range("AN25").Value = 7
$ws.Range("AP25").Value = 27
$ws.Range("AS25").Value = 25
$ws.Range("AU25").Value = 16
$ws.Range("AX25").Value = 23
$ws.Range("BB25").Value = 25
$ws.Range("BC25").Value = 23
$ws.Range("BF25").Value = "2007-12-08"
$ws.Range("D26").Value = 18
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 0.389
$ws.Range("H26").Value = 48.6
$ws.Range("I26").Value = 35.2
$ws.Range("J26").Value = 78.7
$ws.Range("M26").Value = 15.1
$ws.Range("N26").Value = 0.347
$ws.Range("O26").Value = 24.6
$ws.Range("P26").Value = 30.7
$ws.Range("Q26").Value = 0.801
$ws.Range("R26").Value = 10.2
$ws.Range("S26").Value = 29.4
$ws.Range("T26").Value = 39.6
$ws.Range("U26").Value = 16.1
$ws.Range("V26").Value = 15.2
$ws.Range("W26").Value = 7.3
$ws.Range("X26").Value = 3.7
$ws.Range("Y26").Value = 5.4
$ws.Range("Z26").Value = 22.7
$ws.Range("AA26").Value = 25.9
$ws.Range("AB26").Value = 100.2
$ws.Range("AD26").Value = 22
$ws.Range("AF26").Value = 18
$ws.Range("AG26").Value = 20
$ws.Range("AJ26").Value = 24
$ws.Range("AK26").Value = 17
$ws.Range("AN26").Value = 19
$ws.Range("AO26").Value = 1
$ws.Range("AR26").Value = 21
$ws.Range("AV26").Value = 13
$ws.Range("AW26").Value = 16
$ws.Range("AY26").Value = 22
$ws.Range("AZ26").Value = 20
$ws.Range("BF26").Value = "2007-12-08"
$ws.Range("AD27").Value = 2
$ws.Range("AJ27").Value = 22
$ws.Range("AL27").Value = 6
$ws.Range("AO27").Value = 24
$ws.Range("AQ27").Value = 10
$ws.Range("AZ27").Value = 2
$ws.Range("BF27").Value = "2007-12-08"
$ws.Range("AD28").Value = 2
$ws.Range("AH28").Value = 10
$ws.Range("AK28").Value = 24
$ws.Range("AP28").Value = 15
$ws.Range("AQ28").Value = 11
$ws.Range("AX28").Value = 12
$ws.Range("AY28").Value = 20
$ws.Range("BA28").Value = 18
$ws.Range("BF28").Value = "2007-12-08"
$ws.Range("AD29").Value = 2
$ws.Range("AL29").Value = 5
$ws.Range("AR29").Value = 22
$ws.Range("AT29").Value = 24
$ws.Range("AW29").Value = 19
$ws.Range("BF29").Value = "2007-12-08"
$ws.Range("D30").Value = 20
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 0.65
$ws.Range("I30").Value = 40.1
$ws.Range("J30").Value = 80.90000000000001
$ws.Range("K30").Value = 0.496
$ws.Range("M30").Value = 11.3
$ws.Range("N30").Value = 0.356
$ws.Range("Q30").Value = 0.761
$ws.Range("R30").Value = 11.8
$ws.Range("S30").Value = 29.8
$ws.Range("T30").Value = 41.6
$ws.Range("U30").Value = 26.9
$ws.Range("V30").Value = 16.4
$ws.Range("X30").Value = 4.3
$ws.Range("Y30").Value = 5.8
$ws.Range("AA30").Value = 24.1
$ws.Range("AB30").Value = 106.6
$ws.Range("AC30").Value = 7.2
$ws.Range("AD30").Value = 2
$ws.Range("AF30").Value = 6
$ws.Range("AG30").Value = 6
$ws.Range("AJ30").Value = 14
$ws.Range("AK30").Value = 2
$ws.Range("AN30").Value = 13
$ws.Range("AQ30").Value = 13
$ws.Range("AR30").Value = 14
$ws.Range("AS30").Value = 22
$ws.Range("AT30").Value = 19
$ws.Range("AY30").Value = 26
$ws.Range("BB30").Value = 4
$ws.Range("BF30").Value = "2007-12-08"
$ws.Range("AD31").Value = 14
$ws.Range("AH31").Value = 6
$ws.Range("AN31").Value = 21
$ws.Range("AO31").Value = 12
$ws.Range("AP31").Value = 16
$ws.Range("AS31").Value = 14
$ws.Range("AU31").Value = 18
$ws.Range("AW31").Value = 15
$ws.Range("AX31").Value = 5
$ws.Range("AY31").Value = 6
$ws.Range("BA31").Value = 17
$ws.Range("BB31").Value = 11
$ws.Range("BF31").Value = "2007-12-08"
